$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, shifting existing rows 70-114 down to 71-115
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new record
$ws.Cells.Item(70, 1).Value = 11
$ws.Cells.Item(70, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(70, 3).Value = 'Bíobío'
$ws.Cells.Item(70, 4).Value = 44729
$ws.Cells.Item(70, 5).Value = 8
$ws.Cells.Item(70, 6).Value = 100112021
$ws.Cells.Item(70, 7).Value = 'Ají'
$ws.Cells.Item(70, 8).Value = 'Inferno'
$ws.Cells.Item(70, 9).Value = 'Primera'
$ws.Cells.Item(70, 10).Value = 22
$ws.Cells.Item(70, 11).Value = 23000
$ws.Cells.Item(70, 12).Value = 25000
$ws.Cells.Item(70, 13).Value = 23909
$ws.Cells.Item(70, 14).Value = '$/caja 15 kilos'
$ws.Cells.Item(70, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(70, 16).Value = 1594
$ws.Cells.Item(70, 17).Value = 15
$ws.Cells.Item(70, 18).Value = 'Hortaliza'
